$wb = $excel.ActiveWorkbook

# 1. Rename worksheets (same set of typo-swapped names as before, just moved around)
$wb.Worksheets.Item(1).Name = "mEtabolites"
$wb.Worksheets.Item(2).Name = "subjeCtmetabolites"
$wb.Worksheets.Item(3).Name = "subjEctdata"
$wb.Worksheets.Item(4).Name = "vaRmap"
$wb.Worksheets.Item(5).Name = "modEls"
$wb.Worksheets.Item(6).Name = "mOdeloptions"

# 2. Update cell contents on the "modEls" sheet (was "modeLs")
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A1").Value = "modEl"
$ws5.Range("B1").Value = "oUtcomes"
$ws5.Range("C1").Value = "eXposure"
$ws5.Range("D1").Value = "adjusTment"
$ws5.Range("G1").Value = "moDelspec"
$ws5.Range("C2").Value = "cat1 conT1"
$ws5.Range("D3").Value = "cont2 cont3_coRr2"
$ws5.Range("D4").Value = "cat1 cAt1_dup"
$ws5.Range("D5").Value = "cat2 cat3_2"
$ws5.Range("F6").Value = "cont4_mIss > 0"
$ws5.Range("E8").Value = "caT1"

# row 9 on this sheet is cleared out entirely (all values become blank)
$ws5.Range("A9:G9").ClearContents()

# 3. Update cell contents on the "mOdeloptions" sheet (was "modeloptioNs")
$ws6 = $wb.Worksheets.Item(6)

# row 7 becomes blank, and rows 8/9 (the old weights/offset custom option rows) are removed
$ws6.Range("A7:D7").ClearContents()
$ws6.Rows.Item(8).Delete()
$ws6.Rows.Item(8).Delete()
